# GanttChart.xlsx edit: extend the Gantt chart through February/March,
# add a "DEMO" legend swatch next to "PROTOTYPE", adjust the
# Implementation/Development task bars and add a "Final Report" task row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helpers (COM constants, since PowerShell host has no VBA enums here)
# ---------------------------------------------------------------------
# xlPasteFormats = -4122, xlMedium = -4138, xlNone = -4142
# Borders.Item(): 7=xlEdgeLeft 8=xlEdgeTop 9=xlEdgeBottom 10=xlEdgeRight

# ---------------------------------------------------------------------
# 1. New month header blocks: FEBRUARY (V1:Y1, blue) and MARCH (Z1:AC1, red)
#    -- mirrors the existing OCTOBER/NOVEMBER/DECEMBER/JANUARY blocks.
# ---------------------------------------------------------------------
$ws.Range("N1").Copy()
$ws.Range("V1:Y1").PasteSpecial(-4122)
$ws.Range("V1:Y1").Interior.Color = 12611584
$ws.Range("W1").Value = "FEBRUARY"

$ws.Range("N1").Copy()
$ws.Range("Z1:AC1").PasteSpecial(-4122)
$ws.Range("Z1:AC1").Interior.Color = 255
$ws.Range("AA1").Value = "MARCH"

# Week-start date row under the new month headers
$ws.Range("N2").Copy()
$ws.Range("V2:Y2").PasteSpecial(-4122)
$ws.Range("V2:Y2").Interior.Color = 12611584

$ws.Range("N2").Copy()
$ws.Range("Z2:AC2").PasteSpecial(-4122)
$ws.Range("Z2:AC2").Interior.Color = 255

$ws.Range("V2").Value = 44962
$ws.Range("W2").Value = 44969
$ws.Range("X2").Value = 44976
$ws.Range("Y2").Value = 44983
$ws.Range("Z2").Value = 44989
$ws.Range("AA2").Value = 44996
$ws.Range("AB2").Value = 45003
$ws.Range("AC2").Value = 45010

# ---------------------------------------------------------------------
# 2. Legend: give the "PROTOTYPE" swatch box a visible outline, and add
#    a matching "DEMO" swatch box at column R.
# ---------------------------------------------------------------------
$ws.Range("N3:N9").Borders.Item(7).LineStyle = 1
$ws.Range("N3:N9").Borders.Item(7).Weight = -4138
$ws.Range("N3:N9").Borders.Item(10).LineStyle = 1
$ws.Range("N3:N9").Borders.Item(10).Weight = -4138
$ws.Range("N3").Borders.Item(8).LineStyle = 1
$ws.Range("N3").Borders.Item(8).Weight = -4138
$ws.Range("N9").Borders.Item(9).LineStyle = 1
$ws.Range("N9").Borders.Item(9).Weight = -4138

$ws.Range("N3").Copy()
$ws.Range("R3:R9").PasteSpecial(-4122)
$ws.Range("R3").Value = "DEMO"

$ws.Range("R3:R9").Borders.Item(7).LineStyle = 1
$ws.Range("R3:R9").Borders.Item(7).Weight = -4138
$ws.Range("R3:R9").Borders.Item(10).LineStyle = 1
$ws.Range("R3:R9").Borders.Item(10).Weight = -4138
$ws.Range("R3").Borders.Item(8).LineStyle = 1
$ws.Range("R3").Borders.Item(8).Weight = -4138

# ---------------------------------------------------------------------
# 3. Task-bar tweaks on existing rows (3-9): trim the right edge of a
#    couple of boxes that used to be the last column, now that the
#    chart keeps going, and square up the "DEADLINE" marker.
# ---------------------------------------------------------------------
$ws.Range("M4").Borders.Item(10).LineStyle = -4142
$ws.Range("M5").Borders.Item(10).LineStyle = -4142

$ws.Range("N4:N9").Borders.Item(7).LineStyle = 1
$ws.Range("N4:N9").Borders.Item(7).Weight = -4138
$ws.Range("N4:N9").Borders.Item(10).LineStyle = 1
$ws.Range("N4:N9").Borders.Item(10).Weight = -4138
$ws.Range("N9").Borders.Item(9).LineStyle = 1
$ws.Range("N9").Borders.Item(9).Weight = -4138

$ws.Range("O7").Borders.Item(7).LineStyle = -4142
$ws.Range("O7").Borders.Item(9).LineStyle = -4142
$ws.Range("Q7").Borders.Item(10).LineStyle = -4142
$ws.Range("Q7").Borders.Item(9).LineStyle = -4142
$ws.Range("S7").Borders.Item(7).LineStyle = -4142
$ws.Range("S7").Borders.Item(9).LineStyle = -4142
$ws.Range("U7").Borders.Item(9).LineStyle = -4142

# ---------------------------------------------------------------------
# 4. "Implementation" row (8): shrink its bar down to a single week and
#    continue it from O8 all the way out to AC8 (through March).
# ---------------------------------------------------------------------
$ws.Range("J8").Borders.Item(10).LineStyle = 1
$ws.Range("J8").Borders.Item(10).Weight = -4138
$ws.Range("J8").Borders.Item(9).LineStyle = 1
$ws.Range("J8").Borders.Item(9).Weight = -4138

$ws.Range("K8:M8").Interior.Pattern = -4142
$ws.Range("K8:M8").Borders.LineStyle = -4142

$ws.Range("O8").Borders.Item(7).LineStyle = -4142

$ws.Range("O8").Copy()
$ws.Range("V8:AC8").PasteSpecial(-4122)
$ws.Range("AC8").Borders.Item(10).LineStyle = 1
$ws.Range("AC8").Borders.Item(10).Weight = -4138
$ws.Range("U8").Borders.Item(10).LineStyle = -4142

# ---------------------------------------------------------------------
# 5. "Development" row (9): clear its old October/November bar and
#    continue the existing box from O9 out through AC9.
# ---------------------------------------------------------------------
$ws.Range("H9:M9").Interior.Pattern = -4142
$ws.Range("H9:M9").Borders.LineStyle = -4142

$ws.Range("O9").Copy()
$ws.Range("V9:AC9").PasteSpecial(-4122)
$ws.Range("AC9").Borders.Item(10).LineStyle = 1
$ws.Range("AC9").Borders.Item(10).Weight = -4138
$ws.Range("U9").Borders.Item(10).LineStyle = -4142

# ---------------------------------------------------------------------
# 6. New "Final Report" task row (10), with its own bar in March.
# ---------------------------------------------------------------------
$ws.Range("A10").Value = "Final Report"

$ws.Range("O8").Copy()
$ws.Range("Z10:AC10").PasteSpecial(-4122)
$ws.Range("Z10").Borders.Item(7).LineStyle = 1
$ws.Range("Z10").Borders.Item(7).Weight = -4138
$ws.Range("Z10").Borders.Item(8).LineStyle = 1
$ws.Range("Z10").Borders.Item(8).Weight = -4138
$ws.Range("AA10:AC10").Borders.Item(8).LineStyle = 1
$ws.Range("AA10:AC10").Borders.Item(8).Weight = -4138
$ws.Range("AC10").Borders.Item(10).LineStyle = 1
$ws.Range("AC10").Borders.Item(10).Weight = -4138

# ---------------------------------------------------------------------
# 7. Selection cursor, matching where the author left off editing.
# ---------------------------------------------------------------------
$ws.Range("N22").Select()
